$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# Row 4 corresponds to 5a454809-12c8-4e4b-a183-ea6968e68ac6.md
# Column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-19 18:53:56"

# --- zh-cn sheet ---
# Row 4 corresponds to 5a454809-12c8-4e4b-a183-ea6968e68ac6... row
# Column H = "Correspond Handoff Datetime"
# Column K = "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H4").Value = "2016-08-19 18:53:52"
$wsZhCn.Range("K4").Value = "2016-08-19 18:54:16"

# --- de-de sheet ---
# Row 4 corresponds to 5a454809-12c8-4e4b-a183-ea6968e68ac6... row
# Column K = "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("K4").Value = "2016-08-19 18:54:23"
